# Updated symbol list on Mon Dec 12 11:35:42 UTC 2022 with GitHub Actions
#
# Refreshes the scraped coin-price table on Sheet1: updates quoted prices
# (column D) for the rows whose market data shifted since the last run,
# trims a stray "Bestin24h" suffix that leaked into E27, and reorders two
# rows (KickToken / CEJI) whose rank swapped - their whole row (Coin, Link,
# Price, Volume label) moves down/up one slot accordingly.
#
# Values are written with a leading apostrophe and the style is reset back
# to "Normal" right after, so numeric-looking strings such as "283.47"
# persist as literal text (matching how this sheet's cells were originally
# authored as inline strings) instead of being auto-coerced to numbers by
# Excel's normal General-format entry heuristics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'283.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'20.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.207"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06188"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.582"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.487"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Value = "'0.01391"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.1644"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08363"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03667"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.03134"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09130"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.700"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001646"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04682"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006499"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.006194"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001066"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.800"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.3384"
$ws.Range("D25").Style = "Normal"
$ws.Range("E27").Value = "'26AAXTokenAAB"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04720"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "'CEJI"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.008195"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40CEJICEJIBestin24h"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'KickToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.007054"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41KickTokenKICK"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1103"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01149"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006370"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.9994"
$ws.Range("D47").Style = "Normal"
